$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Instagram reel link shared by all three new rows (H.A.D.G tour dates).
$url = "https://www.instagram.com/reel/DGS9b0WM6gg/?igsh=MThjc2lmYWFqZjNpeA=="

# Row 208 ("A208:E208") already carries the plain text/date cell formatting that every
# populated data row in this sheet uses (style "3" for text cells, "4" for the date
# cell). Copying its formats onto the still-empty template rows 209-211 before writing
# values reproduces that look without disturbing styles.xml.
$fmtSrc = $ws.Range("A208:E208")

# ---- Row 209: H.A.D.G @ Club Zimmermanns, Köln, 2025-04-04 -------------------------
$dst209 = $ws.Range("A209:E209")
$fmtSrc.Copy()
$dst209.PasteSpecial(-4122)

$ws.Range("A209").Value = 45751
$ws.Range("B209").Value = "H.A.D.G"
$ws.Range("C209").Value = "Club Zimmermanns"
$ws.Range("D209").Value = "Köln"

$cell209 = $ws.Range("E209")
$cell209.Value = $url
$null = $ws.Hyperlinks.Add($cell209, $url, "", "", $url)

# Give the link text the same underline/blue rich-text run used by every other
# Instagram link in the sheet (Font.ColorIndex 4 renders to the same RGB as this
# workbook's existing indexed-color-11 link runs). Splitting the formatting call
# across the last character keeps it a genuine rich-text run in the shared string
# instead of collapsing into a whole-cell font style.
$chars = $cell209.Characters(1, $url.Length - 1)
$chars.Font.Underline = $true
$chars.Font.ColorIndex = 4
$chars.Font.Name = "Calibri"
$lastChar = $cell209.Characters($url.Length, 1)
$lastChar.Font.Underline = $true
$lastChar.Font.ColorIndex = 4
$lastChar.Font.Name = "Calibri"

# Re-apply the plain row formats so the cell keeps the normal text style (matches
# every other link cell in the sheet, which is styled "3", not Excel's auto
# "Hyperlink" style).
$fmtSrc.Copy()
$dst209.PasteSpecial(-4122)

# ---- Row 210: H.A.D.G @ Die Box, Mönchengladbach, 2025-04-25 -----------------------
$dst210 = $ws.Range("A210:E210")
$fmtSrc.Copy()
$dst210.PasteSpecial(-4122)

$ws.Range("A210").Value = 45772
$ws.Range("B210").Value = "H.A.D.G"
$ws.Range("C210").Value = "Die Box"
$ws.Range("D210").Value = "Mönchengladbach"

# Copy the already-formatted rich-text cell instead of re-building it, so the
# identical link text reuses the same shared-string entry (as in the target file)
# rather than creating a duplicate.
$cell209.Copy($ws.Range("E210"))
$null = $ws.Hyperlinks.Add($ws.Range("E210"), $url, "", "", $url)

$fmtSrc.Copy()
$dst210.PasteSpecial(-4122)

# ---- Row 211: H.A.D.G @ Stollen134, Dortmund, 2025-05-30 ---------------------------
$dst211 = $ws.Range("A211:E211")
$fmtSrc.Copy()
$dst211.PasteSpecial(-4122)

$ws.Range("A211").Value = 45807
$ws.Range("B211").Value = "H.A.D.G"
$ws.Range("C211").Value = "Stollen134"
$ws.Range("D211").Value = "Dortmund"

$cell209.Copy($ws.Range("E211"))
$null = $ws.Hyperlinks.Add($ws.Range("E211"), $url, "", "", $url)

$fmtSrc.Copy()
$dst211.PasteSpecial(-4122)

# Adding hyperlinks registers Excel's built-in "Hyperlink" named cell style even
# though no cell ends up using it (we restore plain formatting above); drop it so
# it doesn't linger unused in the saved styles.
try {
    $wb.Styles.Item("Hyperlink").Delete()
} catch {
}
